$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update simulation result values (marc_s2 simulation data)
$ws.Range("B2").Value = "+108.88 ± 0.00"
$ws.Range("C2").Value = "+108.88 ± 0.00"
$ws.Range("D2").Value = "+23.81 ± 0.02"

$ws.Range("B3").Value = "+108.88 ± 0.00"
$ws.Range("C3").Value = "+108.88 ± 0.00"

$ws.Range("D4").Value = "+21.60 ± 0.02"

$ws.Range("D5").Value = "+2.22 ± 0.01"

$ws.Range("B6").Value = "-108.82 ± 0.01"
$ws.Range("C6").Value = "-108.75 ± 0.01"
$ws.Range("D6").Value = "-23.88 ± 0.02"

$ws.Range("C7").Value = "-21.60 ± 0.02"

$ws.Range("C8").Value = "-2.22 ± 0.01"

$ws.Range("B9").Value = "-13.14 ± 0.01"
$ws.Range("C9").Value = "-0.00 ± 0.00"
$ws.Range("D9").Value = "-13.14 ± 0.01"

$ws.Range("B10").Value = "-0.51 ± 0.00"
$ws.Range("C10").Value = "-0.00 ± 0.00"
$ws.Range("D10").Value = "-0.51 ± 0.00"

$ws.Range("B11").Value = "-90.65 ± 0.02"
$ws.Range("C11").Value = "-81.03 ± 0.02"
$ws.Range("D11").Value = "-9.63 ± 0.02"

$ws.Range("B12").Value = "-4.51 ± 0.00"
$ws.Range("C12").Value = "-3.91 ± 0.00"
$ws.Range("D12").Value = "-0.60 ± 0.00"

$ws.Range("B13").Value = "+1.49 ± 0.00"
$ws.Range("C13").Value = "+1.20 ± 0.00"
$ws.Range("D13").Value = "+0.29 ± 0.00"

$ws.Range("B14").Value = "+5.01 ± 0.01"
$ws.Range("C14").Value = "+4.03 ± 0.01"
$ws.Range("D14").Value = "+4.46 ± 0.02"
